$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.345.42"
$ws.Range("E2").Value = "  +2.41%  "
$ws.Range("D3").Value = "3.170.98"
$ws.Range("E3").Value = "  +1.38%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.52"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.23%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.516"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +6.53%  "
$ws.Range("E9").Value = "  -1.45%  "
$ws.Range("E10").Value = "  +4.12%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.429"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.79%  "
$ws.Range("D12").Value = "3.713.63"
$ws.Range("E12").Value = "  +1.24%  "
$ws.Range("E13").Value = "  +0.40%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.90"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.05%  "
$ws.Range("E15").Value = "  +3.53%  "
$ws.Range("D16").Value = "59.365.78"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").Value = "3.167.05"
$ws.Range("E17").Value = "  +1.11%  "
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.72%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.12"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.28%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "374.43"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.35%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.527"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.76%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.83"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.36%  "
$ws.Range("E25").Value = "  +0.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +15.79%  "
$ws.Range("E27").Value = "  +0.26%  "
$ws.Range("D28").Value = "0.0₃0888"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("E29").Value = "  +1.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("B31").Value = "NEARProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.34"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.80%  "
$ws.Range("B32").Value = "RenderToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.56%  "
$ws.Range("E33").Value = "  -1.13%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.39"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +3.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "156.29"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("E36").Value = "  +3.59%  "
$ws.Range("E37").Value = "  +5.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.45"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.74%  "
$ws.Range("D39").Value = "2.721.04"
$ws.Range("E39").Value = "  +5.96%  "
$ws.Range("E40").Value = "  +0.56%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.26"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.44%  "
$ws.Range("E42").Value = "  +3.96%  "
$ws.Range("E43").Value = "  +2.74%  "
$ws.Range("E44").Value = "  +7.28%  "
$ws.Range("D45").Value = "3.213.02"
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.990"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.52%  "
$ws.Range("E47").Value = "  +0.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0997"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.63%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "20.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.99%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.765"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.22%  "
